$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: "Array" | "Remove Element" | "two pointers, invisible elements;"
$ws.Range("A3").Value = "Array"

# Copy B2's formatting (the highlighted "Notes" style) down into B3, then overwrite its text.
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("B3").Value = "Remove Element"

$ws.Range("C3").Value = "two pointers, invisible elements;"

$ws.Range("C3").Select()
